$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.902.72"
$ws.Range("E2").Value = "  +0.25%  "
$ws.Range("D3").Value = "3.151.36"
$ws.Range("E3").Value = "  +1.36%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "533.35"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.80%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "140.69"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.27%  "
$ws.Range("D8").Value = "3.147.82"
$ws.Range("E8").Value = "  +1.18%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.448"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +2.90%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.24"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.62%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.110"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.39%  "
$ws.Range("E12").Value = "  +4.62%  "
$ws.Range("D13").Value = "3.688.75"
$ws.Range("E13").Value = "  +1.17%  "
$ws.Range("E14").Value = "  +2.96%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "25.69"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -2.03%  "
$ws.Range("E16").Value = "  +0.69%  "
$ws.Range("D17").Value = "57.995.43"
$ws.Range("E17").Value = "  +0.17%  "
$ws.Range("D18").Value = "3.144.72"
$ws.Range("E18").Value = "  +1.09%  "
$ws.Range("E19").Value = "  +0.43%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.86"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.35%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "8.01"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.57%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "356.03"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +5.42%  "
$ws.Range("E23").Value = "  -0.07%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "68.57"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +2.99%  "
$ws.Range("E25").Value = "  +0.73%  "
$ws.Range("E26").Value = "  +1.20%  "
$ws.Range("E27").Value = "  -0.09%  "
$ws.Range("E28").Value = "  +1.57%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.47"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +3.43%  "
$ws.Range("E30").Value = "  +0.02%  "
$ws.Range("E31").Value = "  -2.51%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.91"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +2.50%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "21.30"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +1.78%  "
$ws.Range("E34").Value = "  +0.64%  "
$ws.Range("E35").Value = "  +5.73%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "157.49"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +1.90%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.19"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +1.88%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "26.16"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -3.69%  "
$ws.Range("E39").Value = "  -0.59%  "
$ws.Range("E40").Value = "  +0.98%  "
$ws.Range("E41").Value = "  +11.11%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.10"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +5.09%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.703"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +2.86%  "
$ws.Range("D44").Value = "3.189.01"
$ws.Range("E44").Value = "  +1.09%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "36.73"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.54%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0272"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +5.41%  "
$ws.Range("E47").Value = "  -0.10%  "
$ws.Range("D48").Value = "2.331.71"
$ws.Range("E48").Value = "  +2.21%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.01"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +2.33%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.08"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +1.31%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "20.44"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -1.06%  "
